# This script applies the data corrections and selection changes captured in the
# commit's OOXML diff for Data/Revenue_ton_miles_V2.xlsx.
#
# Content changes:
#   - Sheet "2022": updated weekly Carload/RTM figures for AAR weeks 6-11
#     (columns G:L) across business-unit rows 5-12 (Carloads) and 15-22 (RTMs).
#   - Sheet "2023": filled in previously-zeroed weekly Carload/RTM figures for
#     the same AAR weeks 6-11 (columns G:L) and rows, now that real data is available.
#   - Restored the last-used cell selection on each of those two sheets.
#
# (Workbook-internal fingerprint fields such as fileVersion/revisionPtr GUIDs and the
#  calcFeatures extension list are regenerated by Excel itself on save and are not
#  reachable through the Excel object model, so they are intentionally left alone here.)

$wb = $excel.ActiveWorkbook

$ws2022 = $wb.Worksheets.Item("2022")

# Row 5
$ws2022.Range("G5").Value = 6226
$ws2022.Range("H5").Value = 6020
$ws2022.Range("I5").Value = 5672
$ws2022.Range("J5").Value = 6110
$ws2022.Range("K5").Value = 6443
$ws2022.Range("L5").Value = 6189

# Row 6
$ws2022.Range("G6").Value = 12642
$ws2022.Range("H6").Value = 12720
$ws2022.Range("I6").Value = 12050
$ws2022.Range("J6").Value = 12551
$ws2022.Range("K6").Value = 12215
$ws2022.Range("L6").Value = 13213

# Row 7
$ws2022.Range("G7").Value = 16226
$ws2022.Range("H7").Value = 17051
$ws2022.Range("I7").Value = 15660
$ws2022.Range("J7").Value = 16431
$ws2022.Range("K7").Value = 16066
$ws2022.Range("L7").Value = 17751

# Row 8
$ws2022.Range("G8").Value = 3859
$ws2022.Range("H8").Value = 3467
$ws2022.Range("I8").Value = 3647
$ws2022.Range("J8").Value = 3908
$ws2022.Range("K8").Value = 4421
$ws2022.Range("L8").Value = 4302

# Row 9
$ws2022.Range("G9").Value = 10039
$ws2022.Range("H9").Value = 9397
$ws2022.Range("I9").Value = 8962
$ws2022.Range("J9").Value = 9699
$ws2022.Range("K9").Value = 9906
$ws2022.Range("L9").Value = 10354

# Row 10
$ws2022.Range("G10").Value = 11629
$ws2022.Range("H10").Value = 11309
$ws2022.Range("I10").Value = 11911
$ws2022.Range("J10").Value = 12004
$ws2022.Range("K10").Value = 11457
$ws2022.Range("L10").Value = 11274

# Row 11
$ws2022.Range("G11").Value = 48030
$ws2022.Range("H11").Value = 43862
$ws2022.Range("I11").Value = 43943
$ws2022.Range("J11").Value = 45648
$ws2022.Range("K11").Value = 48738
$ws2022.Range("L11").Value = 48228

# Row 12
$ws2022.Range("G12").Value = 108651
$ws2022.Range("H12").Value = 103826
$ws2022.Range("I12").Value = 101845
$ws2022.Range("J12").Value = 106351
$ws2022.Range("K12").Value = 109246
$ws2022.Range("L12").Value = 111311

# Row 15
$ws2022.Range("G15").Value = 444.64685100000003
$ws2022.Range("H15").Value = 444.742188
$ws2022.Range("I15").Value = 413.33797600000003
$ws2022.Range("J15").Value = 454.99206600000002
$ws2022.Range("K15").Value = 490.41384399999998
$ws2022.Range("L15").Value = 487.29230899999999

# Row 16
$ws2022.Range("G16").Value = 906.57213899999999
$ws2022.Range("H16").Value = 909.75543600000003
$ws2022.Range("I16").Value = 879.67446500000005
$ws2022.Range("J16").Value = 912.27808300000004
$ws2022.Range("K16").Value = 920.15350999999998
$ws2022.Range("L16").Value = 962.79300899999998

# Row 17
$ws2022.Range("G17").Value = 474.906656
$ws2022.Range("H17").Value = 480.75027599999999
$ws2022.Range("I17").Value = 421.68528900000001
$ws2022.Range("J17").Value = 420.55377399999998
$ws2022.Range("K17").Value = 611.49292100000002
$ws2022.Range("L17").Value = 600.02689799999996

# Row 18
$ws2022.Range("G18").Value = 50.502921999999998
$ws2022.Range("H18").Value = 45.840166000000004
$ws2022.Range("I18").Value = 46.436680000000003
$ws2022.Range("J18").Value = 52.560735999999999
$ws2022.Range("K18").Value = 58.026294
$ws2022.Range("L18").Value = 59.000740999999998

# Row 19
$ws2022.Range("G19").Value = 416.71331800000002
$ws2022.Range("H19").Value = 407.08307400000001
$ws2022.Range("I19").Value = 429.830399
$ws2022.Range("J19").Value = 454.415887
$ws2022.Range("K19").Value = 460.30514899999997
$ws2022.Range("L19").Value = 467.89076899999998

# Row 20
$ws2022.Range("G20").Value = 1106.46489
$ws2022.Range("H20").Value = 1042.170805
$ws2022.Range("I20").Value = 1122.6850159999999
$ws2022.Range("J20").Value = 1110.185884
$ws2022.Range("K20").Value = 1006.878017
$ws2022.Range("L20").Value = 1108.8886480000001

# Row 21
$ws2022.Range("G21").Value = 1049.732385
$ws2022.Range("H21").Value = 1003.817105
$ws2022.Range("I21").Value = 1028.636387
$ws2022.Range("J21").Value = 1100.134667
$ws2022.Range("K21").Value = 1108.276386
$ws2022.Range("L21").Value = 1126.500579

# Row 22
$ws2022.Range("G22").Value = 4449.5391609999997
$ws2022.Range("H22").Value = 4334.1590500000002
$ws2022.Range("I22").Value = 4342.286212
$ws2022.Range("J22").Value = 4505.1210970000002
$ws2022.Range("K22").Value = 4655.5461209999994
$ws2022.Range("L22").Value = 4812.3929530000005

$ws2023 = $wb.Worksheets.Item("2023")

# Row 5
$ws2023.Range("G5").Value = 6061
$ws2023.Range("H5").Value = 6454
$ws2023.Range("I5").Value = 6194
$ws2023.Range("J5").Value = 6457
$ws2023.Range("K5").Value = 6539
$ws2023.Range("L5").Value = 6586

# Row 6
$ws2023.Range("G6").Value = 13241
$ws2023.Range("H6").Value = 12782
$ws2023.Range("I6").Value = 12852
$ws2023.Range("J6").Value = 12340
$ws2023.Range("K6").Value = 12515
$ws2023.Range("L6").Value = 11992

# Row 7
$ws2023.Range("G7").Value = 18848
$ws2023.Range("H7").Value = 18497
$ws2023.Range("I7").Value = 17694
$ws2023.Range("J7").Value = 17637
$ws2023.Range("K7").Value = 17353
$ws2023.Range("L7").Value = 16765

# Row 8
$ws2023.Range("G8").Value = 4105
$ws2023.Range("H8").Value = 4163
$ws2023.Range("I8").Value = 4234
$ws2023.Range("J8").Value = 4775
$ws2023.Range("K8").Value = 4288
$ws2023.Range("L8").Value = 4448

# Row 9
$ws2023.Range("G9").Value = 9572
$ws2023.Range("H9").Value = 9461
$ws2023.Range("I9").Value = 10300
$ws2023.Range("J9").Value = 10050
$ws2023.Range("K9").Value = 10278
$ws2023.Range("L9").Value = 10004

# Row 10
$ws2023.Range("G10").Value = 14711
$ws2023.Range("H10").Value = 14796
$ws2023.Range("I10").Value = 12877
$ws2023.Range("J10").Value = 13944
$ws2023.Range("K10").Value = 13830
$ws2023.Range("L10").Value = 14619

# Row 11
$ws2023.Range("G11").Value = 41131
$ws2023.Range("H11").Value = 40590
$ws2023.Range("I11").Value = 39566
$ws2023.Range("J11").Value = 38587
$ws2023.Range("K11").Value = 39019
$ws2023.Range("L11").Value = 38768

# Row 12
$ws2023.Range("G12").Value = 107669
$ws2023.Range("H12").Value = 106743
$ws2023.Range("I12").Value = 103717
$ws2023.Range("J12").Value = 103790
$ws2023.Range("K12").Value = 103822
$ws2023.Range("L12").Value = 103182

# Row 15
$ws2023.Range("G15").Value = 443.31105200000002
$ws2023.Range("H15").Value = 492.69301400000001
$ws2023.Range("I15").Value = 454.65266600000001
$ws2023.Range("J15").Value = 486.92328900000001
$ws2023.Range("K15").Value = 482.55607300000003
$ws2023.Range("L15").Value = 471.20969500000001

# Row 16
$ws2023.Range("G16").Value = 901.21945300000004
$ws2023.Range("H16").Value = 903.55405699999994
$ws2023.Range("I16").Value = 853.31410500000004
$ws2023.Range("J16").Value = 837.315517
$ws2023.Range("K16").Value = 883.31283599999995
$ws2023.Range("L16").Value = 798.79842499999995

# Row 17
$ws2023.Range("G17").Value = 580.84453800000006
$ws2023.Range("H17").Value = 626.280663
$ws2023.Range("I17").Value = 545.30942400000004
$ws2023.Range("J17").Value = 577.22236599999997
$ws2023.Range("K17").Value = 491.91124200000002
$ws2023.Range("L17").Value = 552.191596

# Row 18
$ws2023.Range("G18").Value = 48.154744999999998
$ws2023.Range("H18").Value = 50.694625000000002
$ws2023.Range("I18").Value = 48.305501
$ws2023.Range("J18").Value = 58.714339000000002
$ws2023.Range("K18").Value = 56.807532999999999
$ws2023.Range("L18").Value = 55.489449

# Row 19
$ws2023.Range("G19").Value = 443.45818600000001
$ws2023.Range("H19").Value = 456.44837899999999
$ws2023.Range("I19").Value = 445.14772799999997
$ws2023.Range("J19").Value = 449.23964699999999
$ws2023.Range("K19").Value = 467.123379
$ws2023.Range("L19").Value = 445.75546700000001

# Row 20
$ws2023.Range("G20").Value = 1442.369782
$ws2023.Range("H20").Value = 1398.13257
$ws2023.Range("I20").Value = 1173.676395
$ws2023.Range("J20").Value = 1357.2858630000001
$ws2023.Range("K20").Value = 1353.7928340000001
$ws2023.Range("L20").Value = 1409.8998260000001

# Row 21
$ws2023.Range("G21").Value = 972.01644199999998
$ws2023.Range("H21").Value = 965.34934799999996
$ws2023.Range("I21").Value = 923.63851799999998
$ws2023.Range("J21").Value = 957.78594199999998
$ws2023.Range("K21").Value = 967.64945
$ws2023.Range("L21").Value = 900.66611499999999

# Row 22
$ws2023.Range("G22").Value = 4831.3741980000004
$ws2023.Range("H22").Value = 4893.1526559999993
$ws2023.Range("I22").Value = 4444.0443370000003
$ws2023.Range("J22").Value = 4724.4869630000003
$ws2023.Range("K22").Value = 4703.1533469999995
$ws2023.Range("L22").Value = 4634.0105730000005

# Restore the cell selections recorded in the sheet views.
$ws2022.Activate()
$ws2022.Range("H43").Select()

$ws2023.Activate()
$ws2023.Range("G36").Select()
